# Scheduled runner refresh: update the cached Universalis price/profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) across all eight
# crafter sheets to the latest market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 130.77777
$ws.Range("I5").Value = 130.77777
$ws.Range("K5").Value = 130.77777
$ws.Range("M5").Value = -15.77777
$ws.Range("H18").Value = 1914.3636
$ws.Range("I18").Value = 1914.3636
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1914.3636
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1630.3636
$ws.Range("N18").ClearContents()
$ws.Range("H42").Value = 2223.1667
$ws.Range("I42").Value = 810
$ws.Range("J42").Value = 5049.5
$ws.Range("K42").Value = 2430
$ws.Range("L42").Value = 15148.5
$ws.Range("M42").Value = -2200
$ws.Range("N42").Value = -15608.5
$ws.Range("H45").Value = 4492
$ws.Range("I45").Value = 3499
$ws.Range("J45").Value = 4988.5
$ws.Range("K45").Value = 10497
$ws.Range("L45").Value = 14965.5
$ws.Range("M45").Value = -10305
$ws.Range("N45").Value = -15349.5
$ws.Range("H49").Value = 207.44444
$ws.Range("I49").Value = 151
$ws.Range("J49").Value = 278
$ws.Range("K49").Value = 453
$ws.Range("L49").Value = 834
$ws.Range("M49").Value = -317
$ws.Range("N49").Value = -1106
$ws.Range("H70").Value = 71382.94500000001
$ws.Range("I70").Value = 25497.5
$ws.Range("J70").Value = 77118.625
$ws.Range("K70").Value = 76492.5
$ws.Range("L70").Value = 231355.875
$ws.Range("M70").Value = -76222.5
$ws.Range("N70").Value = -231895.875
$ws.Range("H73").Value = 71382.94500000001
$ws.Range("I73").Value = 25497.5
$ws.Range("J73").Value = 77118.625
$ws.Range("K73").Value = 76492.5
$ws.Range("L73").Value = 231355.875
$ws.Range("M73").Value = -75556.5
$ws.Range("N73").Value = -233227.875
$ws.Range("H74").Value = 4833.364
$ws.Range("I74").Value = 4376.7
$ws.Range("K74").Value = 4376.7
$ws.Range("M74").Value = -3440.7
$ws.Range("H77").Value = 4833.364
$ws.Range("I77").Value = 4376.7
$ws.Range("K77").Value = 21883.5
$ws.Range("M77").Value = -17203.5
$ws.Range("H98").Value = 2017.2826
$ws.Range("I98").Value = 1572.525
$ws.Range("K98").Value = 1572.525
$ws.Range("M98").Value = -74.52500000000009
$ws.Range("H122").Value = 2017.2826
$ws.Range("I122").Value = 1572.525
$ws.Range("K122").Value = 4717.575000000001
$ws.Range("M122").Value = -2267.575000000001
$ws.Range("H132").Value = 2519893.8
$ws.Range("I132").Value = 2761836.8
$ws.Range("K132").Value = 8285510.399999999
$ws.Range("M132").Value = -8282980.399999999
$ws.Range("H137").Value = 25029
$ws.Range("I137").Value = 46778.57
$ws.Range("J137").Value = 5998.125
$ws.Range("K137").Value = 140335.71
$ws.Range("L137").Value = 17994.375
$ws.Range("M137").Value = -137785.71
$ws.Range("N137").Value = -23094.375
$ws.Range("H138").Value = 196384.44
$ws.Range("J138").Value = 635844.75
$ws.Range("L138").Value = 1907534.25
$ws.Range("N138").Value = -1917814.25
$ws.Range("H141").Value = 1905.9
$ws.Range("I141").Value = 1478.2667
$ws.Range("J141").Value = 3188.8
$ws.Range("K141").Value = 4434.800099999999
$ws.Range("L141").Value = 9566.400000000001
$ws.Range("M141").Value = 745.1999000000005
$ws.Range("N141").Value = -19926.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 82510.39999999999
$ws.Range("I32").Value = 71325.56
$ws.Range("K32").Value = 71325.56
$ws.Range("M32").Value = -71038.56
$ws.Range("H92").Value = 1995000
$ws.Range("J92").Value = 3900000
$ws.Range("L92").Value = 3900000
$ws.Range("N92").Value = -3904992
$ws.Range("H132").Value = 1156.9166
$ws.Range("I132").Value = 1048.8
$ws.Range("K132").Value = 3146.4
$ws.Range("M132").Value = -616.3999999999996
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1286.1212
$ws.Range("I86").Value = 1294.5385
$ws.Range("K86").Value = 1294.5385
$ws.Range("M86").Value = -171.5385000000001
$ws.Range("H89").Value = 1286.1212
$ws.Range("I89").Value = 1294.5385
$ws.Range("K89").Value = 6472.692500000001
$ws.Range("M89").Value = -856.692500000001
$ws.Range("H134").Value = 1555.8955
$ws.Range("I134").Value = 1158.3898
$ws.Range("K134").Value = 3475.1694
$ws.Range("M134").Value = -940.1693999999998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11114635
$ws.Range("I31").Value = 33335332
$ws.Range("K31").Value = 33335332
$ws.Range("M31").Value = -33335037
$ws.Range("H34").Value = 11114635
$ws.Range("I34").Value = 33335332
$ws.Range("K34").Value = 33335332
$ws.Range("M34").Value = -33335130
$ws.Range("H58").Value = 2120.9312
$ws.Range("I58").Value = 2089.5
$ws.Range("K58").Value = 2089.5
$ws.Range("M58").Value = -1886.5
$ws.Range("H99").Value = 7435.375
$ws.Range("I99").Value = 6408.4
$ws.Range("J99").Value = 9147
$ws.Range("K99").Value = 6408.4
$ws.Range("L99").Value = 9147
$ws.Range("M99").Value = -4910.4
$ws.Range("N99").Value = -12143
$ws.Range("H105").Value = 1500.6666
$ws.Range("I105").Value = 1174.9166
$ws.Range("K105").Value = 1174.9166
$ws.Range("M105").Value = 572.0834
$ws.Range("H126").Value = 7435.375
$ws.Range("I126").Value = 6408.4
$ws.Range("J126").Value = 9147
$ws.Range("K126").Value = 19225.2
$ws.Range("L126").Value = 27441
$ws.Range("M126").Value = -16755.2
$ws.Range("N126").Value = -32381
$ws.Range("H134").Value = 2559.25
$ws.Range("I134").Value = 2102.4443
$ws.Range("K134").Value = 6307.3329
$ws.Range("M134").Value = -3772.3329
$ws.Range("H135").Value = 92815.414
$ws.Range("J135").Value = 92815.414
$ws.Range("L135").Value = 92815.414
$ws.Range("N135").Value = -102955.414
$ws.Range("H136").Value = 2120.9312
$ws.Range("I136").Value = 2089.5
$ws.Range("K136").Value = 6268.5
$ws.Range("M136").Value = -3718.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4624.125
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 4624.125
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242
$ws.Range("H131").Value = 1531.6522
$ws.Range("J131").Value = 2142.6667
$ws.Range("L131").Value = 6428.000100000001
$ws.Range("N131").Value = -16508.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 614.2917
$ws.Range("I2").Value = 703.5333000000001
$ws.Range("J2").Value = 465.55554
$ws.Range("K2").Value = 703.5333000000001
$ws.Range("L2").Value = 465.55554
$ws.Range("M2").Value = -590.5333000000001
$ws.Range("N2").Value = -691.5555400000001
$ws.Range("H80").Value = 7687.1055
$ws.Range("I80").Value = 4903.5
$ws.Range("K80").Value = 4903.5
$ws.Range("M80").Value = -3905.5
$ws.Range("H83").Value = 7687.1055
$ws.Range("I83").Value = 4903.5
$ws.Range("K83").Value = 24517.5
$ws.Range("M83").Value = -19525.5
$ws.Range("H102").Value = 17230.703
$ws.Range("I102").Value = 18503.734
$ws.Range("K102").Value = 18503.734
$ws.Range("M102").Value = -16881.734
$ws.Range("H126").Value = 3159.25
$ws.Range("I126").Value = 1899.6
$ws.Range("K126").Value = 5698.799999999999
$ws.Range("M126").Value = -3228.799999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H68").Value = 3992.2307
$ws.Range("I68").Value = 3590
$ws.Range("K68").Value = 3590
$ws.Range("M68").Value = -2841
$ws.Range("H71").Value = 3992.2307
$ws.Range("I71").Value = 3590
$ws.Range("K71").Value = 17950
$ws.Range("M71").Value = -14206
$ws.Range("H100").Value = 7065.4
$ws.Range("I100").Value = 4173.5454
$ws.Range("J100").Value = 15018
$ws.Range("K100").Value = 4173.5454
$ws.Range("L100").Value = 15018
$ws.Range("M100").Value = -3632.5454
$ws.Range("N100").Value = -16100
$ws.Range("H122").Value = 5401.1665
$ws.Range("I122").Value = 4281
$ws.Range("K122").Value = 12843
$ws.Range("M122").Value = -10393
$ws.Range("H132").Value = 2348.4243
$ws.Range("I132").Value = 2112
$ws.Range("J132").Value = 4062.5
$ws.Range("K132").Value = 6336
$ws.Range("L132").Value = 12187.5
$ws.Range("M132").Value = -3806
$ws.Range("N132").Value = -17247.5
$ws.Range("H136").Value = 7928.8
$ws.Range("I136").Value = 13222
$ws.Range("K136").Value = 39666
$ws.Range("M136").Value = -37116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2451.127
$ws.Range("I132").Value = 2202.64
$ws.Range("J132").Value = 3406.8462
$ws.Range("K132").Value = 6607.92
$ws.Range("L132").Value = 10220.5386
$ws.Range("M132").Value = -4077.92
$ws.Range("N132").Value = -15280.5386
$ws.Range("H136").Value = 33299.168
$ws.Range("I136").Value = 47678
$ws.Range("J136").Value = 4541.5
$ws.Range("K136").Value = 143034
$ws.Range("L136").Value = 13624.5
$ws.Range("M136").Value = -140484
$ws.Range("N136").Value = -18724.5
